# Update the 3rd data row (row 3) of the active sheet from the "Mini"
# sample record to a new "Suresh" record, then move the active selection
# to J4 (mirrors the author's manual edits captured in the workbook diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters here: it controls the order new entries are appended to
# the shared string table, matching the target workbook's string order
# (male, Suresh, suresh@hihi.hi, active).
$ws.Range("B3").Value = "male"
$ws.Range("A3").Value = "Suresh"
$ws.Range("C3").Value = "suresh@hihi.hi"
$ws.Range("D3").Value = "active"

# Move/record the current selection as it was left in the saved workbook.
$ws.Range("J4").Select()
